# Correlation Summary.xlsx - data exploration update
# Adds "Weight_Features" / PROD_FREQ based method-comparison block (cols G:K)
# to the "SOLO MEN" sheet and removes the old scratch notes in C32:C34.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SOLO MEN")

# ---------------------------------------------------------------------
# 1) Remove the old leftover notes in C32:C34 ("[7,8,9]" / "[7,8,10]" / "[7,8,11]")
# ---------------------------------------------------------------------
$ws.Range("C32:C34").ClearContents()

# ---------------------------------------------------------------------
# 2) New header note above the table
# ---------------------------------------------------------------------
$ws.Range("G4").Value = "CAMBIO DE MÉTODO… AHORA PESAMOS CADA FEATURE CON BASE EN PROD_FREQ"

# ---------------------------------------------------------------------
# 3) New "divisiones" comparison blocks in columns G:K
#    Each block: a "Método N divisiones…" label, then a merged G:H cell
#    labelled "HDC ('Prod_Freq', 30, MEN)" with 3 numeric results in I:K.
# ---------------------------------------------------------------------

function Add-DivisionesBlock {
    param(
        [string]$LabelCell,
        [string]$LabelText,
        [int]$Row,
        [double]$I,
        [double]$J,
        [double]$K,
        [bool]$HasK = $true
    )

    $ws.Range($LabelCell).Value = $LabelText

    $hdrRange = "G" + $Row + ":H" + $Row
    $ws.Range("B19").Copy()
    $ws.Range($hdrRange).PasteSpecial(-4122)
    $ws.Application.CutCopyMode = $false
    $ws.Range($hdrRange).Merge()
    $ws.Range("G" + $Row).Value = "HDC ('Prod_Freq', 30, MEN)"

    $ws.Range("C8").Copy()
    $ws.Range("I" + $Row).PasteSpecial(-4122)
    $ws.Range("J" + $Row).PasteSpecial(-4122)
    $ws.Range("K" + $Row).PasteSpecial(-4122)
    $ws.Application.CutCopyMode = $false

    $ws.Range("I" + $Row).Value = $I
    $ws.Range("J" + $Row).Value = $J
    if ($HasK) {
        $ws.Range("K" + $Row).Value = $K
    }
}

# Método 5 divisiones… (row 12 label, row 13 data)
Add-DivisionesBlock "G12" "Método 5 divisiones…" 13 0.6523 0.6546 0.6448 $true

# Método 4 divisiones… (row 15 label, row 16 data)
Add-DivisionesBlock "G15" "Método 4 divisiones…" 16 0.6526 0.6565 0.6574 $true

# Método 3 divisiones… (row 18 label, row 19 data)
Add-DivisionesBlock "G18" "Método 3 divisiones…" 19 0.6552 0.65801 0.6674 $true

# Método 2 divisiones… (row 21 label, row 22 data -- no K value)
Add-DivisionesBlock "G21" "Método 2 divisiones…" 22 0.6644 0.666 0 $false

# Método 6 divisiones… (row 9 label, row 10 data) -- added last
Add-DivisionesBlock "G9" "Método 6 divisiones…" 10 0.6424 0.6389 0.6432 $true

# ---------------------------------------------------------------------
# 4) Sheet view: drop the old scroll/selection pointing at the removed notes
# ---------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("L11").Select()
